$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '28.897.35'
Set-TextValue $ws.Range('E2') '  +3.02%  '
Set-TextValue $ws.Range('D3') '1.884.29'
Set-TextValue $ws.Range('E3') '  +3.07%  '
Set-TextValue $ws.Range('D4') '1.005'
Set-TextValue $ws.Range('E4') '  +0.54%  '
Set-TextValue $ws.Range('D5') '324.67'
Set-TextValue $ws.Range('E5') '  -1.28%  '
Set-TextValue $ws.Range('D6') '1.005'
Set-TextValue $ws.Range('E6') '  +0.46%  '
Set-TextValue $ws.Range('D7') '0.4675'
Set-TextValue $ws.Range('D8') '0.3940'
Set-TextValue $ws.Range('E8') '  +1.80%  '
Set-TextValue $ws.Range('D9') '0.07930'
Set-TextValue $ws.Range('E9') '  +0.59%  '
Set-TextValue $ws.Range('D10') '0.9841'
Set-TextValue $ws.Range('E10') '  +2.41%  '
Set-TextValue $ws.Range('D11') '22.39'
Set-TextValue $ws.Range('E11') '  +2.13%  '
Set-TextValue $ws.Range('D12') '1.866.79'
Set-TextValue $ws.Range('E12') '  +2.68%  '
Set-TextValue $ws.Range('D13') '5.758'
Set-TextValue $ws.Range('E13') '  +1.77%  '
Set-TextValue $ws.Range('D14') '7.019'
Set-TextValue $ws.Range('E14') '  +1.80%  '
Set-TextValue $ws.Range('D15') '0.06991'
Set-TextValue $ws.Range('E15') '  +2.19%  '
Set-TextValue $ws.Range('D16') '88.91'
Set-TextValue $ws.Range('E16') '  +2.72%  '
Set-TextValue $ws.Range('E17') '  +0.57%  '
Set-TextValue $ws.Range('E18') '  +1.34%  '
Set-TextValue $ws.Range('D19') '17.00'
Set-TextValue $ws.Range('E19') '  +2.04%  '
Set-TextValue $ws.Range('E20') '  +0.27%  '
Set-TextValue $ws.Range('D21') '28.890.76'
Set-TextValue $ws.Range('E21') '  +2.97%  '
Set-TextValue $ws.Range('D22') '5.357'
Set-TextValue $ws.Range('E22') '  +0.57%  '
Set-TextValue $ws.Range('D23') '11.11'
Set-TextValue $ws.Range('E23') '  +0.95%  '
Set-TextValue $ws.Range('D24') '2.123'
Set-TextValue $ws.Range('E24') '  +1.41%  '
Set-TextValue $ws.Range('D25') '2.093.02'
Set-TextValue $ws.Range('E25') '  +2.81%  '
Set-TextValue $ws.Range('D26') '153.83'
Set-TextValue $ws.Range('E26') '  +1.08%  '
Set-TextValue $ws.Range('D27') '19.43'
Set-TextValue $ws.Range('E27') '  +1.11%  '
Set-TextValue $ws.Range('D28') '5.799'
Set-TextValue $ws.Range('E28') '  +0.77%  '
Set-TextValue $ws.Range('D29') '2.005'
Set-TextValue $ws.Range('E29') '  +1.77%  '
Set-TextValue $ws.Range('D30') '119.99'
Set-TextValue $ws.Range('E30') '  +2.64%  '
Set-TextValue $ws.Range('E31') '  +1.81%  '
Set-TextValue $ws.Range('E32') '  +0.26%  '
Set-TextValue $ws.Range('D33') '5.324'
Set-TextValue $ws.Range('E33') '  +0.62%  '
Set-TextValue $ws.Range('E34') '  +3.27%  '
Set-TextValue $ws.Range('D35') '3.347'
Set-TextValue $ws.Range('E35') '  +0.09%  '
Set-TextValue $ws.Range('E36') '  -0.16%  '
Set-TextValue $ws.Range('D37') '0.02127'
Set-TextValue $ws.Range('E37') '  -0.73%  '
Set-TextValue $ws.Range('E38') '  +1.36%  '
Set-TextValue $ws.Range('D39') '7.899'
Set-TextValue $ws.Range('E39') '  +3.66%  '
Set-TextValue $ws.Range('D40') '0.5736'
Set-TextValue $ws.Range('E40') '  +2.61%  '
Set-TextValue $ws.Range('D41') '0.1803'
Set-TextValue $ws.Range('E41') '  +2.03%  '
Set-TextValue $ws.Range('D42') '10.02'
Set-TextValue $ws.Range('E42') '  +0.95%  '
Set-TextValue $ws.Range('D43') '0.07309'
Set-TextValue $ws.Range('E43') '  +4.33%  '
Set-TextValue $ws.Range('D44') '11.87'
Set-TextValue $ws.Range('E44') '  +2.32%  '
Set-TextValue $ws.Range('D45') '0.5356'
Set-TextValue $ws.Range('E45') '  +1.84%  '
Set-TextValue $ws.Range('D46') '1.173'
Set-TextValue $ws.Range('E46') '  -4.36%  '
Set-TextValue $ws.Range('D47') '2.130'
Set-TextValue $ws.Range('E47') '  -4.42%  '
Set-TextValue $ws.Range('D48') '1.851'
Set-TextValue $ws.Range('E48') '  +1.34%  '
Set-TextValue $ws.Range('D49') '114.28'
Set-TextValue $ws.Range('E49') '  +2.58%  '
Set-TextValue $ws.Range('D50') '2.377'
Set-TextValue $ws.Range('E50') '  +3.02%  '
Set-TextValue $ws.Range('D51') '1.005'
Set-TextValue $ws.Range('E51') '  +0.45%  '
